$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.390.89"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "1.638.94"
$ws.Range("E3").Value = "  -1.73%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'211.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.71%  "
$ws.Range("E6").Value = "  +3.82%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'23.05"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.71%  "
$ws.Range("E9").Value = "  -2.31%  "
$ws.Range("E10").Value = "  -2.15%  "
$ws.Range("D11").Value = "'0.0890"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.18%  "
$ws.Range("D12").Value = "1.870.48"
$ws.Range("D13").Value = "1.651.21"
$ws.Range("E13").Value = "  -0.98%  "
$ws.Range("D14").Value = "'4.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.76%  "
$ws.Range("E15").Value = "  -0.64%  "
$ws.Range("D16").Value = "'64.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.94%  "
$ws.Range("D17").Value = "27.362.42"
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("D18").Value = "'229.44"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.80%  "
$ws.Range("D19").Value = "0.0₃0720"
$ws.Range("E19").Value = "  -1.29%  "
$ws.Range("D20").Value = "'7.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.71%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("E22").Value = "  -3.58%  "
$ws.Range("D23").Value = "'9.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.28%  "
$ws.Range("E24").Value = "  -0.75%  "
$ws.Range("D25").Value = "'147.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("D26").Value = "'6.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.25%  "
$ws.Range("E27").Value = "  +1.52%  "
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("D29").Value = "'15.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.61%  "
$ws.Range("E30").Value = "  -4.22%  "
$ws.Range("E31").Value = "  -3.25%  "
$ws.Range("E32").Value = "  -2.30%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").Value = "1.408.21"
$ws.Range("E34").Value = "  -3.94%  "
$ws.Range("E35").Value = "  +0.27%  "
$ws.Range("E36").Value = "  -0.26%  "
$ws.Range("E37").Value = "  -1.92%  "
$ws.Range("E38").Value = "  -5.46%  "
$ws.Range("E39").Value = "  -3.46%  "
$ws.Range("E40").Value = "  +0.89%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("E42").Value = "  -1.56%  "
$ws.Range("D43").Value = "'5.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.01%  "
$ws.Range("E44").Value = "  +0.66%  "
$ws.Range("D45").Value = "'0.792"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.51%  "
$ws.Range("D46").Value = "'64.48"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.20%  "
$ws.Range("D47").Value = "1.780.18"
$ws.Range("E47").Value = "  -1.62%  "
$ws.Range("E48").Value = "  -4.65%  "
$ws.Range("D49").Value = "'87.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.32%  "
$ws.Range("D50").Value = "0.0₆0105"
$ws.Range("E50").Value = "  -2.04%  "
$ws.Range("E51").Value = "  -3.67%  "
